$wb = $excel.ActiveWorkbook

# --- Resumen sheet: update Maximo (C2) ---
$resumen = $wb.Worksheets.Item("Resumen")
$resumen.Range("C2").Value = 701.8706887231054

# --- Solucion sheet: update Pedido/Salida assignments ---
$solucion = $wb.Worksheets.Item("Solucion")
$solucion.Range("A2").Value = "Pedido_48"
$solucion.Range("B2").Value = "S001"
$solucion.Range("A3").Value = "Pedido_27"
$solucion.Range("B3").Value = "S025"
$solucion.Range("A4").Value = "Pedido_32"
$solucion.Range("B4").Value = "S041"
$solucion.Range("A5").Value = "Pedido_15"
$solucion.Range("B5").Value = "S065"
$solucion.Range("A6").Value = "Pedido_75"
$solucion.Range("B6").Value = "S005"
$solucion.Range("A7").Value = "Pedido_38"
$solucion.Range("B7").Value = "S029"
$solucion.Range("A8").Value = "Pedido_35"
$solucion.Range("B8").Value = "S069"
$solucion.Range("A9").Value = "Pedido_6"
$solucion.Range("B9").Value = "S045"
$solucion.Range("A10").Value = "Pedido_7"
$solucion.Range("B10").Value = "S002"
$solucion.Range("A11").Value = "Pedido_63"
$solucion.Range("B11").Value = "S026"
$solucion.Range("A12").Value = "Pedido_3"
$solucion.Range("B12").Value = "S042"
$solucion.Range("A13").Value = "Pedido_13"
$solucion.Range("B13").Value = "S030"
$solucion.Range("A14").Value = "Pedido_28"
$solucion.Range("B14").Value = "S066"
$solucion.Range("A15").Value = "Pedido_77"
$solucion.Range("B15").Value = "S006"
$solucion.Range("A16").Value = "Pedido_47"
$solucion.Range("B16").Value = "S003"
$solucion.Range("A17").Value = "Pedido_36"
$solucion.Range("B17").Value = "S027"
$solucion.Range("A18").Value = "Pedido_24"
$solucion.Range("B18").Value = "S070"
$solucion.Range("A19").Value = "Pedido_51"
$solucion.Range("B19").Value = "S046"
$solucion.Range("A20").Value = "Pedido_23"
$solucion.Range("B20").Value = "S007"
$solucion.Range("A21").Value = "Pedido_5"
$solucion.Range("B21").Value = "S067"
$solucion.Range("A22").Value = "Pedido_52"
$solucion.Range("B22").Value = "S031"
$solucion.Range("A23").Value = "Pedido_44"
$solucion.Range("B23").Value = "S043"
$solucion.Range("A24").Value = "Pedido_57"
$solucion.Range("B24").Value = "S004"
$solucion.Range("A25").Value = "Pedido_43"
$solucion.Range("B25").Value = "S047"
$solucion.Range("A26").Value = "Pedido_8"
$solucion.Range("B26").Value = "S008"
$solucion.Range("A27").Value = "Pedido_19"
$solucion.Range("B27").Value = "S028"
$solucion.Range("A28").Value = "Pedido_39"
$solucion.Range("B28").Value = "S071"
$solucion.Range("A29").Value = "Pedido_55"
$solucion.Range("B29").Value = "S044"
$solucion.Range("A30").Value = "Pedido_76"
$solucion.Range("B30").Value = "S032"
$solucion.Range("A31").Value = "Pedido_46"
$solucion.Range("B31").Value = "S009"
$solucion.Range("A32").Value = "Pedido_54"
$solucion.Range("B32").Value = "S068"
$solucion.Range("A33").Value = "Pedido_37"
$solucion.Range("B33").Value = "S048"
$solucion.Range("A34").Value = "Pedido_74"
$solucion.Range("B34").Value = "S033"
$solucion.Range("A35").Value = "Pedido_31"
$solucion.Range("B35").Value = "S072"
$solucion.Range("A36").Value = "Pedido_22"
$solucion.Range("B36").Value = "S013"
$solucion.Range("A37").Value = "Pedido_18"
$solucion.Range("B37").Value = "S049"
$solucion.Range("A38").Value = "Pedido_34"
$solucion.Range("B38").Value = "S073"
$solucion.Range("A39").Value = "Pedido_79"
$solucion.Range("B39").Value = "S010"
$solucion.Range("A40").Value = "Pedido_29"
$solucion.Range("B40").Value = "S037"
$solucion.Range("A41").Value = "Pedido_49"
$solucion.Range("B41").Value = "S053"
$solucion.Range("A42").Value = "Pedido_78"
$solucion.Range("B42").Value = "S014"
$solucion.Range("A43").Value = "Pedido_62"
$solucion.Range("B43").Value = "S034"
$solucion.Range("A44").Value = "Pedido_25"
$solucion.Range("B44").Value = "S077"
$solucion.Range("A45").Value = "Pedido_69"
$solucion.Range("B45").Value = "S011"
$solucion.Range("A46").Value = "Pedido_40"
$solucion.Range("B46").Value = "S050"
$solucion.Range("A47").Value = "Pedido_33"
$solucion.Range("B47").Value = "S038"
$solucion.Range("A48").Value = "Pedido_72"
$solucion.Range("B48").Value = "S074"
$solucion.Range("A49").Value = "Pedido_30"
$solucion.Range("B49").Value = "S015"
$solucion.Range("A50").Value = "Pedido_1"
$solucion.Range("B50").Value = "S078"
$solucion.Range("A51").Value = "Pedido_59"
$solucion.Range("B51").Value = "S054"
$solucion.Range("A52").Value = "Pedido_53"
$solucion.Range("B52").Value = "S012"
$solucion.Range("A53").Value = "Pedido_45"
$solucion.Range("B53").Value = "S035"
$solucion.Range("A54").Value = "Pedido_21"
$solucion.Range("B54").Value = "S075"
$solucion.Range("A55").Value = "Pedido_11"
$solucion.Range("B55").Value = "S051"
$solucion.Range("A56").Value = "Pedido_50"
$solucion.Range("B56").Value = "S016"
$solucion.Range("A57").Value = "Pedido_10"
$solucion.Range("B57").Value = "S079"
$solucion.Range("A58").Value = "Pedido_17"
$solucion.Range("B58").Value = "S039"
$solucion.Range("A59").Value = "Pedido_12"
$solucion.Range("B59").Value = "S017"
$solucion.Range("A60").Value = "Pedido_73"
$solucion.Range("B60").Value = "S055"
$solucion.Range("A61").Value = "Pedido_2"
$solucion.Range("B61").Value = "S036"
$solucion.Range("A62").Value = "Pedido_80"
$solucion.Range("B62").Value = "S076"
$solucion.Range("A63").Value = "Pedido_70"
$solucion.Range("B63").Value = "S052"
$solucion.Range("A64").Value = "Pedido_67"
$solucion.Range("B64").Value = "S021"
$solucion.Range("A65").Value = "Pedido_14"
$solucion.Range("B65").Value = "S056"
$solucion.Range("A66").Value = "Pedido_4"
$solucion.Range("B66").Value = "S080"
$solucion.Range("A67").Value = "Pedido_26"
$solucion.Range("B67").Value = "S040"
$solucion.Range("A68").Value = "Pedido_20"
$solucion.Range("B68").Value = "S018"
$solucion.Range("A69").Value = "Pedido_41"
$solucion.Range("B69").Value = "S057"
$solucion.Range("A70").Value = "Pedido_58"
$solucion.Range("B70").Value = "S061"
$solucion.Range("A71").Value = "Pedido_66"
$solucion.Range("B71").Value = "S022"
$solucion.Range("A72").Value = "Pedido_64"
$solucion.Range("B72").Value = "S058"
$solucion.Range("A73").Value = "Pedido_65"
$solucion.Range("B73").Value = "S019"
$solucion.Range("A74").Value = "Pedido_56"
$solucion.Range("B74").Value = "S023"
$solucion.Range("A75").Value = "Pedido_16"
$solucion.Range("B75").Value = "S062"
$solucion.Range("A76").Value = "Pedido_9"
$solucion.Range("B76").Value = "S059"
$solucion.Range("A77").Value = "Pedido_61"
$solucion.Range("B77").Value = "S020"
$solucion.Range("A78").Value = "Pedido_60"
$solucion.Range("B78").Value = "S063"
$solucion.Range("A79").Value = "Pedido_71"
$solucion.Range("B79").Value = "S024"
$solucion.Range("A80").Value = "Pedido_68"
$solucion.Range("B80").Value = "S060"
$solucion.Range("A81").Value = "Pedido_42"
$solucion.Range("B81").Value = "S064"

# --- Metricas sheet: update Tiempo values ---
$metricas = $wb.Worksheets.Item("Metricas")
$metricas.Range("B2").Value = 688.2933560384907
$metricas.Range("B3").Value = 476.3639042058601
$metricas.Range("B4").Value = 701.8706887231054
$metricas.Range("B5").Value = 465.7289977294841
